# "Natmi following Dr Hou advice"
# Re-run of the Bmp2->Acvr2b NATMI edge table with the EC (endothelial
# cell) cluster added as a third sending/target cluster alongside FAPs
# and sCs, which changes every downstream specificity/weight figure and
# grows the table from 6 to 9 combination rows (3x3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the expression / specificity metrics.
$data = @(
    @("ECs", "Bmp2", "Acvr2b", "ECs", 2, 0.6666666666666666, 9.163165666666666, 27.489497, 0.2800251397703982, 0.2800251397703982, 3, 1, 0.5292956666666667, 1.587887, 0.135651968140022, 0.1356519681400219, 4.850023880315445, 43.650214922839, 0.03798596133853924, 0.03798596133853924),
    @("ECs", "Bmp2", "Acvr2b", "FAPs", 2, 0.6666666666666666, 9.163165666666666, 27.489497, 0.2800251397703982, 0.2800251397703982, 3, 1, 1.362890666666667, 4.088672000000001, 0.3492921120199358, 0.3492921120199358, 12.48839296422045, 112.395536677984, 0.09781057248908011, 0.09781057248908012),
    @("ECs", "Bmp2", "Acvr2b", "sCs", 2, 0.6666666666666666, 9.163165666666666, 27.489497, 0.2800251397703982, 0.2800251397703982, 3, 1, 2.009678666666666, 6.029036, 0.5150559198400423, 0.5150559198400423, 18.41501855943244, 165.735167034892, 0.1442286059427788, 0.1442286059427789),
    @("FAPs", "Bmp2", "Acvr2b", "ECs", 3, 1, 12.06704733333333, 36.201142, 0.3687673822623249, 0.3687673822623249, 3, 1, 0.5292956666666667, 1.587887, 0.135651968140022, 0.1356519681400219, 6.38703586299489, 57.483322766954, 0.05002402118972819, 0.05002402118972818),
    @("FAPs", "Bmp2", "Acvr2b", "FAPs", 3, 1, 12.06704733333333, 36.201142, 0.3687673822623249, 0.3687673822623249, 3, 1, 1.362890666666667, 4.088672000000001, 0.3492921120199358, 0.3492921120199358, 16.44606618482489, 148.014595663424, 0.1288075377944705, 0.1288075377944705),
    @("FAPs", "Bmp2", "Acvr2b", "sCs", 3, 1, 12.06704733333333, 36.201142, 0.3687673822623249, 0.3687673822623249, 3, 1, 2.009678666666666, 6.029036, 0.5150559198400423, 0.5150559198400423, 24.25088759545688, 218.257988359112, 0.1899358232781262, 0.1899358232781262),
    @("sCs", "Bmp2", "Acvr2b", "ECs", 3, 1, 11.49244066666667, 34.477322, 0.3512074779672769, 0.351207477967277, 3, 1, 0.5292956666666667, 1.587887, 0.135651968140022, 0.1356519681400219, 6.082899044290445, 54.74609139861401, 0.04764198561175451, 0.04764198561175451),
    @("sCs", "Bmp2", "Acvr2b", "FAPs", 3, 1, 11.49244066666667, 34.477322, 0.3512074779672769, 0.351207477967277, 3, 1, 1.362890666666667, 4.088672000000001, 0.3492921120199358, 0.3492921120199358, 15.66294012182045, 140.966461096384, 0.1226740017363852, 0.1226740017363852),
    @("sCs", "Bmp2", "Acvr2b", "sCs", 3, 1, 11.49244066666667, 34.477322, 0.3512074779672769, 0.351207477967277, 3, 1, 2.009678666666666, 6.029036, 0.5150559198400423, 0.5150559198400423, 23.09611283573244, 207.865015521592, 0.1808914906191372, 0.1808914906191372)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
